$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Widen column D slightly (raw OOXML width 7 -> 8).
# Excel's ColumnWidth property (in characters) is offset from the raw
# stored width by the default column padding, so 7.15 characters round-trips
# to a stored width of exactly 8.
$ws.Columns.Item(4).ColumnWidth = 7.15

# Update "Lao" (column D) values
$ws.Range("D2").Value = 1580
$ws.Range("D3").Value = 10

# Update "Total" (column G) values to match
$ws.Range("G2").Value = 1580
$ws.Range("G3").Value = 10
